$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain plain text (matches the source
# data which already stores these as text, e.g. "58.129.50", "  -1.48%  ").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '58.198.09'
$ws.Range("E2").Value = '  -0.81%  '

$ws.Range("D3").Value = '2.478.51'
$ws.Range("E3").Value = '  -0.47%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").Value = '519.70'
$ws.Range("E5").Value = '  -2.75%  '

$ws.Range("D6").Value = '131.80'
$ws.Range("E6").Value = '  -3.26%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -1.03%  '

$ws.Range("D9").Value = '0.0995'
$ws.Range("E9").Value = '  -0.93%  '

$ws.Range("D10").Value = '0.157'
$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("D11").Value = '5.36'
$ws.Range("E11").Value = '  +0.91%  '

$ws.Range("D12").Value = '0.344'
$ws.Range("E12").Value = '  -0.41%  '

$ws.Range("D13").Value = '2.917.51'
$ws.Range("E13").Value = '  -1.32%  '

$ws.Range("D14").Value = '58.126.29'
$ws.Range("E14").Value = '  -0.94%  '

$ws.Range("D15").Value = '22.13'
$ws.Range("E15").Value = '  -3.08%  '

$ws.Range("E16").Value = '  -1.28%  '

$ws.Range("D17").Value = '2.482.76'
$ws.Range("E17").Value = '  -0.87%  '

$ws.Range("D18").Value = '10.86'
$ws.Range("E18").Value = '  -1.69%  '

$ws.Range("D19").Value = '4.19'
$ws.Range("E19").Value = '  -1.78%  '

$ws.Range("D20").Value = '320.45'
$ws.Range("E20").Value = '  -0.56%  '

$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '5.75'
$ws.Range("E22").Value = '  -2.15%  '

$ws.Range("D23").Value = '64.29'
$ws.Range("E23").Value = '  -1.53%  '

$ws.Range("D24").Value = '0.410'
$ws.Range("E24").Value = '  -2.15%  '

$ws.Range("D25").Value = '0.998'
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("E26").Value = '  -2.73%  '

$ws.Range("E27").Value = '  -2.11%  '

$ws.Range("D28").Value = '0.0₃0758'
$ws.Range("E28").Value = '  -0.61%  '

$ws.Range("D29").Value = '1.71'
$ws.Range("E29").Value = '  -3.12%  '

$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").Value = '6.34'
$ws.Range("E30").Value = '  -5.15%  '

$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '166.54'
$ws.Range("E31").Value = '  -0.30%  '

$ws.Range("E32").Value = '  +1.07%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("E34").Value = '  +0.25%  '

$ws.Range("D35").Value = '18.13'
$ws.Range("E35").Value = '  -1.30%  '

$ws.Range("E36").Value = '  -9.48%  '

$ws.Range("D37").Value = '3.98'
$ws.Range("E37").Value = '  -2.07%  '

$ws.Range("E38").Value = '  -3.00%  '

$ws.Range("D39").Value = '0.793'
$ws.Range("E39").Value = '  -2.04%  '

$ws.Range("D40").Value = '3.49'
$ws.Range("E40").Value = '  -2.60%  '

$ws.Range("D41").Value = '276.47'
$ws.Range("E41").Value = '  -2.79%  '

$ws.Range("D42").Value = '5.04'
$ws.Range("E42").Value = '  -3.15%  '

$ws.Range("E43").Value = '  -1.09%  '

$ws.Range("D44").Value = '126.69'
$ws.Range("E44").Value = '  -3.43%  '

$ws.Range("D45").Value = '0.0907'
$ws.Range("E45").Value = '  -1.63%  '

$ws.Range("D46").Value = '0.0490'
$ws.Range("E46").Value = '  -2.99%  '

$ws.Range("E47").Value = '  -2.19%  '

$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("D49").Value = '1.743.51'
$ws.Range("E49").Value = '  -1.21%  '

$ws.Range("E50").Value = '  -1.34%  '

$ws.Range("D51").Value = '4.69'
$ws.Range("E51").Value = '  -1.11%  '
